# Update Krankenversicherungsbeitraege.xlsx
#
# Renames two rows on "Tabelle1":
#   A5: "Beitragsbemessungsgrenze GKV Ost"  -> "Beitragsbemessungsgrenze GKV"
#   A6: "Beitragsbemessungsgrenze GKV West" -> "Jahresarbeitsentgeltgrenze GKV"
# and moves the active selection to A9 (matches the saved sheetView state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A5").Value = "Beitragsbemessungsgrenze GKV"
$ws.Range("A6").Value = "Jahresarbeitsentgeltgrenze GKV"

$ws.Range("A9").Select() | Out-Null
